$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently sits right after
#    the H1 title ("Play Fruit Shop Free: Fun & easy online slot game").
# ---------------------------------------------------------------------------
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. Before the final paragraph (the italic "Create a cartoon-style image..."
#    AI-image prompt paragraph), insert a bold "Play Fruit Shop Free: Fun &
#    easy online slot game" paragraph, and replace the final paragraph's text
#    with the meta-description copy (keeping its italic formatting).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$replacementXml = "<w:p $wNs>" +
    "<w:r/>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruit Shop Free: Fun &amp; easy online slot game</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
    "<w:r/>" +
    "<w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Fruit Shop and play for free. Enjoy straightforward gameplay, winning opportunities, and a high RTP rate of 96.70%.</w:t></w:r>" +
    "</w:p>"

$null = $lastPara.Range.InsertXML($replacementXml)
